$d = $word.ActiveDocument

# --- The results table (descr() output) that had a placeholder label ---
$t = $d.Tables(2)

# Table overall width: 10947 -> 9507 dxa (twips); COM works in points (1 pt = 20 dxa)
$t.PreferredWidth = 475.35

# First column width: 3232 -> 1792 dxa
$t.Columns(1).Width = 89.6

# Row 3 (the "My custom variable label" placeholder row): trHeight 454 -> 457 dxa
$t.Rows(3).Height = 22.85

# That row's first cell: indent left 400 -> 100 dxa, make the run bold, and
# swap the placeholder text for the real variable name "Sepal.Length".
$cell = $t.Cell(3, 1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="left"/><w:spacing w:after="40" w:before="40"/><w:ind w:firstLine="0" w:left="100" w:right="100"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:eastAsia="Arial" w:cs="Arial"/><w:b/><w:sz w:val="22"/><w:szCs w:val="22"/><w:color w:val="111111"/></w:rPr><w:t xml:space="preserve">Sepal.Length</w:t></w:r></w:p>'
$cell.Range.Paragraphs(1).Range.InsertXML($xml)
